# Fixed problem in publications:
# - The "venue" cell for the first publication row (C2) contained the wrong
#   text ("ECCV UNCV Workshop 2022 - Extended Abstract"); it should just be
#   "ECCV UNCV Workshop".
# - Update the active cell selection on the "publications" sheet to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("publications")

$ws.Range("C2").Value = "ECCV UNCV Workshop"

$ws.Activate() | Out-Null
$ws.Range("C2").Select() | Out-Null
